$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 32, shifting existing rows 32:51 down to 33:52
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new data record
$ws.Cells.Item(32, 1).Value = 11
$ws.Cells.Item(32, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(32, 3).Value = "Bíobío"
$ws.Cells.Item(32, 4).Value = 45062
$ws.Cells.Item(32, 5).Value = 8
$ws.Cells.Item(32, 6).Value = 100112026
$ws.Cells.Item(32, 7).Value = "Haba"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 100
$ws.Cells.Item(32, 11).Value = 20000
$ws.Cells.Item(32, 12).Value = 21000
$ws.Cells.Item(32, 13).Value = 20500
$ws.Cells.Item(32, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(32, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(32, 16).Value = 820
$ws.Cells.Item(32, 17).Value = 25
$ws.Cells.Item(32, 18).Value = "Hortaliza"
